$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle existing ToParse (column H) flags
$ws.Range("H2").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 1

# Add a new row to the table (Table1) so ref/autoFilter/dimension grow with it
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()

# New row 29 - Silpo / Brovary
$ws.Range("A29").Value = "Сільпо"
$ws.Range("B29").Value = "https://silpo.ua/"
$ws.Range("C29").Value = "Торти, тістечка"
$ws.Range("D29").Value = "https://silpo.ua/category/torty-tistechka-663"
$ws.Range("E29").Value = "silpo_parser"
$ws.Range("F29").Value = "Бровари"
$ws.Range("G29").Value = "вул. Київська, 156"
$ws.Range("H29").Value = 1

# Match the hyperlink-style look used by the other rows' B/D cells (no actual hyperlink)
$ws.Range("B29").Style = $ws.Range("B28").Style
$ws.Range("D29").Style = $ws.Range("D28").Style

# Extend the whole-number data validation on column H down through the new row
$ws.Range("H2:H29").Validation.Delete()
$ws.Range("H2:H29").Validation.Add(1, 1, 1, 0, 1)

$ws.Range("H16").Select()
